# Auto-generated PowerShell Excel COM-interop script
# Applies per-cell text updates to sheet1 (cryptos list) as described in the diff.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Range('D2').Value = '61.150.82'
$ws.Range('E2').Value = '  +0.20%  '
$ws.Range('D3').Value = '2.920.53'
$ws.Range('E3').Value = '  +0.12%  '
$ws.Range('E4').Value = '  +0.05%  '
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D5').Value = '589.78'
$ws.Range('E5').Value = '  +0.60%  '
$ws.Range('E6').Value = '  -1.56%  '
$ws.Range('E7').Value = '  +0.10%  '
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D8').Value = '0.507'
$ws.Range('E8').Value = '  +0.94%  '
$ws.Range('D9').Value = '2.917.83'
$ws.Range('E9').Value = '  +0.01%  '
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D10').Value = '6.81'
$ws.Range('E10').Value = '  +0.95%  '
$ws.Range('E11').Value = '  -0.05%  '
$ws.Range('E12').Value = '  -1.09%  '
$ws.Range('E13').Value = '  +1.00%  '
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D14').Value = '33.81'
$ws.Range('E14').Value = '  -2.09%  '
$ws.Range('D16').Value = '3.404.50'
$ws.Range('E16').Value = '  +0.16%  '
$ws.Range('D17').Value = '61.108.16'
$ws.Range('E17').Value = '  +0.21%  '
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D18').Value = '6.71'
$ws.Range('E18').Value = '  -1.95%  '
$ws.Range('D19').Value = '2.929.38'
$ws.Range('E19').Value = '  +0.49%  '
$ws.Range('D20').NumberFormat = "@"
$ws.Range('D20').Value = '431.12'
$ws.Range('E20').Value = '  +1.07%  '
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D21').Value = '13.48'
$ws.Range('E21').Value = '  -1.43%  '
$ws.Range('E22').Value = '  +1.78%  '
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D23').Value = '7.09'
$ws.Range('E23').Value = '  -1.31%  '
$ws.Range('E24').Value = '  -0.04%  '
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D25').Value = '11.02'
$ws.Range('E25').Value = '  -0.17%  '
$ws.Range('D26').NumberFormat = "@"
$ws.Range('D26').Value = '2.24'
$ws.Range('E26').Value = '  +2.61%  '
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D27').Value = '12.04'
$ws.Range('E27').Value = '  +1.94%  '
$ws.Range('E28').Value = '  +0.02%  '
$ws.Range('E29').Value = '  +6.76%  '
$ws.Range('E30').Value = '  +0.09%  '
$ws.Range('E31').Value = '  +0.17%  '
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D32').Value = '7.16'
$ws.Range('E32').Value = '  -1.74%  '
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D33').Value = '26.50'
$ws.Range('E33').Value = '  -1.08%  '
$ws.Range('E34').Value = '  +1.33%  '
$ws.Range('D35').Value = '0.0₃0866'
$ws.Range('E35').Value = '  +2.27%  '
$ws.Range('E36').Value = '  +0.34%  '
$ws.Range('E37').Value = '  -0.67%  '
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D38').Value = '3.09'
$ws.Range('E38').Value = '  +2.32%  '
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D39').Value = '49.94'
$ws.Range('E39').Value = '  +0.15%  '
$ws.Range('E40').Value = '  -1.86%  '
$ws.Range('E41').Value = '  -0.95%  '
$ws.Range('E42').Value = '  -1.98%  '
$ws.Range('E43').Value = '  +0.07%  '
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D44').Value = '39.55'
$ws.Range('E44').Value = '  -5.46%  '
$ws.Range('D45').NumberFormat = "@"
$ws.Range('D45').Value = '377.62'
$ws.Range('E45').Value = '  +0.23%  '
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D46').Value = '0.0348'
$ws.Range('E46').Value = '  +0.66%  '
$ws.Range('D47').Value = '2.709.50'
$ws.Range('E47').Value = '  +1.98%  '
$ws.Range('D48').NumberFormat = "@"
$ws.Range('D48').Value = '131.70'
$ws.Range('E48').Value = '  -1.18%  '
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D50').Value = '24.21'
$ws.Range('E50').Value = '  -5.21%  '
$ws.Range('E51').Value = '  +0.05%  '
